# Updates the cryptocurrency listing cells (Coin / Link / Price / Volume(1h))
# to the latest scraped values, matching a GitHub Actions refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new literal text to write. The source workbook
# stores every cell (including the numeric-looking "Price" column) as text,
# so every value here is written verbatim as a string rather than a number.
$updates = [ordered]@{
    "D2" = '63.391.52'
    "E2" = '  -7.12%  '
    "D3" = '3.534.73'
    "E3" = '  -2.40%  '
    "E4" = '  -0.01%  '
    "D5" = '389.48'
    "E5" = '  -6.64%  '
    "D6" = '121.73'
    "E6" = '  -6.35%  '
    "D7" = '3.525.91'
    "E7" = '  -2.30%  '
    "D8" = '0.582'
    "E8" = '  -11.42%  '
    "D9" = '1.00'
    "E9" = '  +0.25%  '
    "E10" = '  -11.37%  '
    "E11" = '  -22.59%  '
    "D12" = '0.0000322'
    "E12" = '  -24.06%  '
    "D13" = '38.35'
    "E13" = '  -8.68%  '
    "D14" = '4.087.41'
    "E14" = '  -2.24%  '
    "E15" = '  -7.59%  '
    "E16" = '  -3.05%  '
    "D17" = '3.514.14'
    "E17" = '  -3.28%  '
    "E18" = '  +2.65%  '
    "D19" = '18.57'
    "E19" = '  -7.28%  '
    "D20" = '63.420.39'
    "E20" = '  -6.85%  '
    "D21" = '1.01'
    "E21" = '  -9.48%  '
    "D22" = '388.45'
    "E22" = '  -15.57%  '
    "D23" = '13.81'
    "E23" = '  +4.22%  '
    "D24" = '80.43'
    "E24" = '  -9.86%  '
    "E25" = '  -6.44%  '
    "E26" = '  +10.44%  '
    "D27" = '33.43'
    "E27" = '  -5.78%  '
    "E28" = '  -9.28%  '
    "E29" = '  -14.91%  '
    "D30" = '11.67'
    "E30" = '  -4.23%  '
    "D31" = '2.66'
    "E31" = '  -4.54%  '
    "E32" = '  -7.17%  '
    "E33" = '  -7.40%  '
    "E34" = '  -5.68%  '
    "E35" = '  +0.08%  '
    "D36" = '36.31'
    "E36" = '  -9.83%  '
    "D37" = '53.31'
    "E37" = '  -4.78%  '
    "E38" = '  -10.78%  '
    "E39" = '  -0.02%  '
    "E40" = '  +3.61%  '
    "E41" = '  -12.40%  '
    "B42" = 'EnergySwap'
    "C42" = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    "D42" = '25.80'
    "E42" = '  +21.96%  '
    "B43" = 'ApeXProtocol'
    "C43" = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
    "D43" = '3.02'
    "E43" = '  +15.13%  '
    "B44" = 'Monero'
    "C44" = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    "D44" = '139.37'
    "E44" = '  -6.38%  '
    "D45" = '0.0₃0606'
    "E45" = '  -23.66%  '
    "D46" = '1.95'
    "E46" = '  +0.72%  '
    "D47" = '4.07'
    "E47" = '  -4.11%  '
    "B48" = 'WEMIXToken'
    "C48" = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    "D48" = '2.47'
    "E48" = '  -9.17%  '
    "B49" = 'LidoDAOToken'
    "C49" = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    "D49" = '3.03'
    "E49" = '  -6.11%  '
    "E50" = '  -10.72%  '
    "E51" = '  -9.76%  '
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $cell = $ws.Range($addr)

    # "Price" values such as "1.00", "121.73" or "0.0000322" are valid
    # numeric literals; left alone, Excel's COM layer auto-converts them to
    # numbers on assignment (dropping the formatted text, e.g. "1.00" -> 1).
    # Values with multiple separators (e.g. "63.391.52") or other
    # non-numeric characters are never coerced, so only force the text
    # format for the ambiguous, purely-numeric-looking strings.
    if ($value -match '^-?\d+(\.\d+)?$') {
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}
